# Set formatted_headnote method in child class for exporting.
#
# The document's final "SectionHeadnote" paragraph currently holds one big
# run of concatenated text that actually represents a whole section's worth
# of resources (a case, another case, and the start of "Section Two").
# Break it up into its proper paragraphs/styles, matching the structure
# used for "Section One" / its resources elsewhere in the casebook export.

$d = $word.ActiveDocument
$cr = [char]13
$lf = [char]10

# Locate the paragraph holding the old combined blob of text.
$find = $d.Content
$find.Find.Execute( `
    "1.1Case of the District Number 1This is the body of case 1.1.2Case of the District Number 2highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $find.Find.Found) {
    throw "could not locate target SectionHeadnote paragraph"
}

$target = $find.Paragraphs(1)

# Build the replacement text: 12 paragraphs, separated by CR (paragraph
# marks). NOTE: don't mix literal LF characters into this same assignment -
# when a CR is present anywhere in a Range.Text value, embedded LFs also
# turn into paragraph breaks; we add the (literal, non-breaking) trailing
# newlines afterwards with their own Range.Text assignments instead.
$parts = @(
    "What is a corporation?",
    "1.1",
    "Case of the District Number 1",
    "",
    "This is the body of case 1.",
    "1.2",
    "Case of the District Number 2",
    "This is an annotatable resource in the casebook.",
    "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;",
    "2",
    "Section Two",
    "This is the second chapter of the casebook."
)

$newText = [string]::Join($cr, $parts)
$target.Range.Text = $newText

# Walk backwards from the now-last paragraph (the new final "SectionHeadnote")
# applying the correct style to each of the 12 new paragraphs.
$styles = @(
    "SectionHeadnote",
    "SectionTitle",
    "SectionNumber",
    "CaseText",
    "ResourceHeadnote",
    "ResourceTitle",
    "ResourceNumber",
    "CaseText",
    "ResourceHeadnote",
    "ResourceTitle",
    "ResourceNumber",
    "SectionHeadnote"
)

$p = $d.Paragraphs.Last
foreach ($styleName in $styles) {
    $p.Style = $styleName
    $prev = $p.Previous()
    $p = $prev
}

# Three of the runs end with a literal trailing newline character (matching
# the original document's embedded "\n" between blobs) - add those back as
# plain text, not paragraph marks.
function Set-TrailingNewline($searchText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r.Find.Found) {
        throw "could not find text for trailing newline: $searchText"
    }
    $r.Text = $searchText + $lf
}

Set-TrailingNewline "This is an annotatable resource in the casebook."
Set-TrailingNewline "highlighted: content to highlight; elided: content to elide; replaced: content to replace; commented: content to comment; highlighted2: second highlight content;"
Set-TrailingNewline "This is the second chapter of the casebook."

# Wrap the resource/section numbers in bookmarks, matching the anchors
# referenced by the table of contents hyperlinks.
function Add-NumberBookmark($searchText, $bookmarkName) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r.Find.Found) {
        throw "could not find text to bookmark: $searchText"
    }
    $d.Bookmarks.Add($bookmarkName, $r)
}

Add-NumberBookmark "1.1" "_auto_toc_2"
Add-NumberBookmark "1.2" "_auto_toc_3"
Add-NumberBookmark "2" "_auto_toc_4"
